# Fruta / hortaliza, semanal
# Re-order the weekly records in rows 3-7 (cyclic shuffle of the data rows),
# keeping columns A, B, C, E, F, G, H, I, J, Q, R, T unchanged since those
# are identical across the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, K, L, M, N, O, P, S for rows 3..7
$data = @{
    3 = @{ D = 44313; K = "Mankaki"; L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; S = 1194 }
    4 = @{ D = 44305; K = "Mankaki"; L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; S = 1361 }
    5 = @{ D = 44301; K = "Hachiya"; L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; S = 1139 }
    6 = @{ D = 44342; K = "Mankaki"; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; S = 1361 }
    7 = @{ D = 44699; K = "Mankaki"; L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; S = 1639 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
